$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07559333333333333
$ws.Range("H2").Value = 0.22678
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 3.844784628635555
$ws.Range("R2").Value = 34.60306165772
$ws.Range("S2").Value = 0.3434314568613803
$ws.Range("T2").Value = 0.3434314568613804

# Row 3
$ws.Range("G3").Value = 0.07559333333333333
$ws.Range("H3").Value = 0.22678
$ws.Range("M3").Value = 43.683024
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("Q3").Value = 3.30214539424
$ws.Range("R3").Value = 29.71930854816
$ws.Range("S3").Value = 0.294960761928139
$ws.Range("T3").Value = 0.294960761928139

# Row 4
$ws.Range("G4").Value = 0.07559333333333333
$ws.Range("H4").Value = 0.22678
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 2.770012297695555
$ws.Range("R4").Value = 24.93011067926
$ws.Range("S4").Value = 0.2474285170192034
$ws.Range("T4").Value = 0.2474285170192035

# Row 5
$ws.Range("G5").Value = 0.07559333333333333
$ws.Range("H5").Value = 0.22678
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 1.278259958722222
$ws.Range("R5").Value = 11.5043396285
$ws.Range("S5").Value = 0.1141792641912772
$ws.Range("T5").Value = 0.1141792641912772
